$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Collapse the "Crud - json-server - bootstrap ui ..." title, the
#    "Npm init", "Npm i json-server", "Mettre dans le script de ...",
#    "Créer le fichier db.json ...", "Npm start ...", "Bootstrap ui"
#    and "Npm install bootstrap ..." paragraphs down to their plain
#    (proofErr-free) run layout. Re-emitting each paragraph's own
#    WordOpenXML back into itself normalizes away the w:proofErr
#    spell-check bookmarks and merges runs that only existed because
#    of those bookmarks, while paragraphs whose separate runs matter
#    for other reasons (diverging run-level formatting/rsid) keep
#    their original split.
# ------------------------------------------------------------------
$paraIds = @(
    "2B3ACE71", # Crud - json-server - bootsptrap ui - angular reactive form
    "639F11AD", # Npm init
    "19613BB2", # Npm i json-server
    "0777ED1F", # Mettre dans le script de package.json : "start": "..."
    "30DD0D48", # Créer le fichier db.json dans le dossier ...
    "124C8239", # Npm start pour lancer le server json-server
    "5485683A", # Bootstrap ui
    "716FDF0F"  # Npm install bootstrap,  il fait aller dans le dossier du projet
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    $xml = $r.WordOpenXML
    foreach ($pid in $paraIds) {
        if ($xml -match ('w14:paraId="' + $pid + '"')) {
            $r.InsertXML($xml)
            break
        }
    }
}

# ------------------------------------------------------------------
# 2) "  Json server" paragraph: the leading two-space run keeps its
#    own (differently-rsid'd) run, only "Json" / " server" merge into
#    a single proofErr-free run.
# ------------------------------------------------------------------
$jsonServerXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4614A9DB" w14:textId="5BD84BDC" w:rsidR="007319CC" w:rsidRDefault="00CC1FB9" w:rsidP="007959E1"><w:pPr><w:pStyle w:val="MonTitre2"/></w:pPr><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:t>Json server</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Json server") {
        $p.Range.InsertXML($jsonServerXml)
        break
    }
}

# ------------------------------------------------------------------
# 3) Replace the first trailing blank paragraph (right after the
#    "Npm install bootstrap..." bullet) with a "bonjou" paragraph
#    indented 360 twips (18 pt) from the left.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.WordOpenXML -match 'w14:paraId="3D298AEA"') {
        $p.Range.Text = "bonjou"
        $p.Format.LeftIndent = 18
        break
    }
}

Write-Output "done"
